$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stats")

# Add new row of data (curchar) below existing rows
$ws.Range("A5").Value = 99
$ws.Range("B5").Value = "curchar"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# Match the selection left behind in the saved file
$ws.Range("G5").Select()
